# feat: Filter out example rows from Excel templates during upload
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column width adjustments ---
# Note: the stored OOXML <col width> is derived from ColumnWidth via a
# characters->pixels->characters round trip (6px/char + 5px padding), so we
# dial in the COM-visible ColumnWidth that lands exactly on the target
# stored width (target - 5/6) rather than the target integer itself.
$ws.Columns.Item(1).ColumnWidth = (9 - 5/6)    # A: 11 -> 9
$ws.Columns.Item(4).ColumnWidth = (19 - 5/6)   # D: 20 -> 19
$ws.Columns.Item(5).ColumnWidth = (6 - 5/6)    # E: 7 -> 6
$ws.Columns.Item(13).ColumnWidth = (14 - 5/6)  # M: 15 -> 14
$ws.Columns.Item(14).ColumnWidth = (6 - 5/6)   # N: 9 -> 6

# --- Row 2: turn into an "example" row with placeholder data ---
$ws.Range("A2").Value = "예시-0001"
$ws.Range("B2").Value = "거래처명 예시"
$ws.Range("H2").Value = "09:00"
$ws.Range("I2").Value = "17:00"
$ws.Range("L2").Value = "홍길동"
$ws.Range("N2").Value = ""

# --- Row 3: clear out the second sample customer, keep only timing/flag data ---
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "09:00"
$ws.Range("G3").Value = "17:00"
$ws.Range("H3").Value = "09:00"
$ws.Range("I3").Value = "17:00"
$ws.Range("K3").Value = 30
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = ""
